$d = $word.ActiveDocument

# Locate the "GIS & Geospatial Analysis Consulting" paragraph (under the
# PARTNER - Siege Analytics heading) and insert three new bullet paragraphs
# directly after it, before the existing "Lead comprehensive research..." bullet.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("GIS & Geospatial Analysis Consulting", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)

if ($found) {
    $anchorPara = $searchRange.Paragraphs(1)
    $currentIndex = $anchorPara.Index

    $bullets = @(
        "• Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels",
        "• Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide",
        "• Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis"
    )

    foreach ($bulletText in $bullets) {
        $currentPara = $d.Paragraphs($currentIndex)
        $currentPara.Range.InsertParagraphAfter()
        $currentIndex = $currentIndex + 1
        $newPara = $d.Paragraphs($currentIndex)
        $newPara.Range.Text = $bulletText
    }
}
